$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.001.40'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '3.786.68'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '''599.03'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '''169.80'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('D7').Value = '3.784.01'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '''0.524'
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('E10').Value = '  -2.49%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = '''0.0000262'
$ws.Range('E13').Value = '  -5.45%  '
$ws.Range('D14').Value = '''36.78'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = '4.423.59'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').Value = '3.786.67'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '69.002.18'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').Value = '''18.15'
$ws.Range('E18').Value = '  -2.76%  '
$ws.Range('D19').Value = '''7.08'
$ws.Range('E19').Value = '  -1.90%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '''11.00'
$ws.Range('E21').Value = '  +4.15%  '
$ws.Range('D22').Value = '''471.10'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('D24').Value = '''84.72'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = '''0.0000147'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('D26').Value = '''2.23'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '''12.18'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').Value = '''10.23'
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = '3.937.12'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('D31').Value = '''2.81'
$ws.Range('E31').Value = '  -3.33%  '
$ws.Range('D32').Value = '''7.45'
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''30.27'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''2.23'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').Value = '''9.38'
$ws.Range('E35').Value = '  +2.45%  '
$ws.Range('D37').Value = '3.743.05'
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('D39').Value = '''3.52'
$ws.Range('E39').Value = '  -9.27%  '
$ws.Range('E40').Value = '  +1.48%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').Value = '''5.87'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('D47').Value = '''43.93'
$ws.Range('E47').Value = '  +11.87%  '
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').Value = '''46.12'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('D50').Value = '''399.82'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').Value = '''146.24'
$ws.Range('E51').Value = '  +4.59%  '
